$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.471.33'
$ws.Range('E2').Value = '  +0.54%  '

$ws.Range('D3').Value = '1.943.12'
$ws.Range('E3').Value = '  -2.11%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.56%  '

$ws.Range('E6').Value = '  -2.83%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.28'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -7.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.367'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.03%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '55.83'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.21%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0839'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.68%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.104'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.33%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.826'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.67%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.49'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.63%  '

$ws.Range('D15').Value = '2.228.70'
$ws.Range('E15').Value = '  -2.21%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.61'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.43%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.26'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.76%  '

$ws.Range('D18').Value = '1.933.78'
$ws.Range('E18').Value = '  -3.05%  '

$ws.Range('D19').Value = '36.407.64'
$ws.Range('E19').Value = '  +0.60%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.20%  '

$ws.Range('E21').Value = '  -0.40%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.39%  '

$ws.Range('E24').Value = '  +0.00%  '

$ws.Range('E25').Value = '  -1.83%  '

$ws.Range('E26').Value = '  -0.96%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.04%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.58'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.81%  '

$ws.Range('E29').Value = '  -2.53%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.124'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -10.15%  '

$ws.Range('E31').Value = '  -1.97%  '

$ws.Range('E32').Value = '  +0.01%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.69'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.71%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0632'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.37%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.78%  '

$ws.Range('E36').Value = '  -2.03%  '

$ws.Range('E37').Value = '  -0.11%  '

$ws.Range('E38').Value = '  -3.14%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.15'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.01%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.41%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0971'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.50%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.96'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.21%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.18'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.90%  '

$ws.Range('E44').Value = '  -2.26%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.34%  '

$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.352.14'
$ws.Range('E46').Value = '  -0.11%  '

$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.04'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.91%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.77'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.96%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.21'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.51%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.81'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.91%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.22%  '
